$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 18
$ws.Range("B3").Value = "Update index.py"
$ws.Range("C3").Value = "riya-morankar"
$ws.Range("D3").Value = "N/A"
$ws.Range("E3").Value = "from edit1 to main"

# Force the date to be stored as literal text "2025-06-17" rather than
# being auto-converted into a date serial number, then strip the
# resulting explicit cell style so it matches the other text cells.
$ws.Range("F3").Value = "'2025-06-17"
$ws.Range("F3").Style = "Normal"
